# "added 4wk low sales check"
# Updates the Forecast Comparison sheet to reflect the new "Low Volume Season"
# trend/low-sales logic, and refreshes the dependent Summary totals.

$wb = $excel.ActiveWorkbook
$fc = $wb.Worksheets.Item("Forecast Comparison")
$sm = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet -------------------------------------------

# Row 2 (W10)
$fc.Range("D2").Value = 6
$fc.Range("G2").Value = "Low Volume Season"
$fc.Range("H2").Value = 25.83
$fc.Range("L2").Value = 0.9399999999999999

# Row 3 (W11)
$fc.Range("D3").Value = 1
$fc.Range("G3").Value = "Low Volume Season"
$fc.Range("H3").Value = 149
$fc.Range("L3").Value = 1.03

# Row 4 (W12)
$fc.Range("D4").Value = 0
$fc.Range("G4").Value = "Low Volume Season"
$fc.Range("H4").Value = ""
$fc.Range("L4").Value = 1.09

# Row 5 (W13)
$fc.Range("D5").Value = 0
$fc.Range("G5").Value = "Low Volume Season"
$fc.Range("H5").Value = ""
$fc.Range("L5").Value = 1.12

# Row 6 (W14)
$fc.Range("D6").Value = 0
$fc.Range("G6").Value = "Low Volume Season"
$fc.Range("H6").Value = ""
$fc.Range("L6").Value = 1.1

# Row 7 (W15)
$fc.Range("D7").Value = 0
$fc.Range("G7").Value = "Low Volume Season"
$fc.Range("H7").Value = ""
$fc.Range("L7").Value = 1.05

# Row 8 (W16)
$fc.Range("D8").Value = 0
$fc.Range("G8").Value = "Low Volume Season"
$fc.Range("H8").Value = ""
$fc.Range("L8").Value = 0.91

# Row 9 (W17)
$fc.Range("D9").Value = 0
$fc.Range("G9").Value = "Low Volume Season"
$fc.Range("H9").Value = ""
$fc.Range("L9").Value = 0.96

# Row 10 (W18)
$fc.Range("D10").Value = 0
$fc.Range("G10").Value = "Low Volume Season"
$fc.Range("H10").Value = ""
$fc.Range("L10").Value = 1.18

# Row 11 (W19)
$fc.Range("D11").Value = 0
$fc.Range("G11").Value = "Low Volume Season"
$fc.Range("H11").Value = ""
$fc.Range("L11").Value = 1.07

# Row 12 (W20)
$fc.Range("D12").Value = 0
$fc.Range("G12").Value = "Low Volume Season"
$fc.Range("H12").Value = ""
$fc.Range("L12").Value = 0.95

# Row 13 (W21)
$fc.Range("D13").Value = 0
$fc.Range("G13").Value = "Low Volume Season"
$fc.Range("H13").Value = ""
$fc.Range("L13").Value = 0.95

# Row 14 (W22)
$fc.Range("D14").Value = 0
$fc.Range("G14").Value = "Low Volume Season"
$fc.Range("H14").Value = ""
$fc.Range("L14").Value = 1.12

# Row 15 (W23)
$fc.Range("D15").Value = 0
$fc.Range("G15").Value = "Low Volume Season"
$fc.Range("H15").Value = ""
$fc.Range("L15").Value = 1.01

# Row 16 (W24)
$fc.Range("D16").Value = 0
$fc.Range("G16").Value = "Low Volume Season"
$fc.Range("H16").Value = ""
$fc.Range("I16").Value = "Low"
$fc.Range("J16").Value = "Normal"
$fc.Range("L16").Value = 1.03

# Row 17 (W25)
$fc.Range("D17").Value = 0
$fc.Range("G17").Value = "Low Volume Season"
$fc.Range("H17").Value = ""
$fc.Range("I17").Value = "Low"
$fc.Range("J17").Value = "Normal"
$fc.Range("L17").Value = 0.83

# --- Summary sheet ----------------------------------------------------------

$sm.Range("B9").Value = "7"    # Total Forecast (16 Weeks)
$sm.Range("B10").Value = "7"   # Total Forecast (8 Weeks)
$sm.Range("B11").Value = "7"   # Total Forecast (4 Weeks)
$sm.Range("B12").Value = "6"   # Max Forecast
$sm.Range("B14").Value = "0"   # Min Forecast

Write-Host "Applied 4wk low sales check updates"
